# Add a new worksheet "Result1" after the existing "Result" sheet, mirroring the
# Result sheet layout but re-ordering the Low52W/DiffFrom52WHigh columns and adding
# a new DiffFrom52WLow column (LTP / Low52W * 100).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Result1"

$headers = @("Symbol", "LTP", "High52W", "DiffFrom52WHigh", "DiffFrom52WLow")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Each row: Symbol, LTP, High52W, DiffFrom52WHigh, Low52W, DiffFrom52WLow
$rows = @()
$rows += @(,@("OFSS", 7698.0, 13220.0, 58.22995461422088, 7038.0, 109.37766410912192))
$rows += @(,@("DIXON", 12061.0, 18700.0, 64.49732620320856, 11646.0, 103.56345526360982))
$rows += @(,@("COLPAL", 2077.1, 2975.0, 69.81848739495798, 2051.0, 101.27254997562163))
$rows += @(,@("SUPREMEIND", 3366.7, 4817.95, 69.8782677279756, 3095.0, 108.77867528271405))
$rows += @(,@("JUBLFOOD", 558.25, 796.75, 70.06589268904926, 548.4, 101.79613420860686))
$rows += @(,@("GODREJPROP", 2006.0, 2848.0, 70.43539325842697, 1900.0, 105.57894736842105))
$rows += @(,@("TIINDIA", 2612.8, 3693.6, 70.73857483214209, 2407.1, 108.54555273981141))
$rows += @(,@("SUZLON", 52.7, 74.3, 70.92866756393002, 46.15, 114.19284940411703))
$rows += @(,@("PAGEIND", 36075.0, 50590.0, 71.30855900375569, 35310.0, 102.16652506372132))
$rows += @(,@("OBEROIRLTY", 1668.6, 2337.7, 71.3778500235274, 1451.95, 114.92131271738006))
$rows += @(,@("MANKIND", 2190.9, 2998.4, 73.06897011739595, 2090.0, 104.82775119617226))
$rows += @(,@("PIIND", 3235.0, 4330.0, 74.71131639722864, 2951.1, 109.62014164209955))
$rows += @(,@("INDUSINDBK", 865.5, 1086.55, 79.65579126593346, 606.0, 142.82178217821783))
$rows += @(,@("FORTIS", 884.75, 1104.3, 80.11862718464187, 577.0, 153.3362218370884))
$rows += @(,@("POLICYBZR", 1830.0, 2246.9, 81.44554719836218, 1311.35, 139.55084454951006))
$rows += @(,@("IRCTC", 685.15, 831.75, 82.37451157198676, 656.0, 104.44359756097562))
$rows += @(,@("COFORGE", 1664.9, 2005.36, 83.02249970080186, 1194.01, 139.43769315164866))
$rows += @(,@("SBICARD", 862.0, 1027.25, 83.91336091506449, 663.05, 130.00527863660358))
$rows += @(,@("NHPC", 79.2, 92.34, 85.76998050682262, 71.0, 111.5492957746479))
$rows += @(,@("OIL", 425.0, 494.55, 85.9367101405318, 325.0, 130.76923076923077))
$rows += @(,@("BSE", 2637.2, 3030.0, 87.03630363036302, 1227.33, 214.87293555930353))
$rows += @(,@("DABUR", 502.35, 577.0, 87.06239168110919, 433.3, 115.93584121855527))
$rows += @(,@("AUROPHARMA", 1187.4, 1356.2, 87.55345819200707, 1010.0, 117.56435643564357))
$rows += @(,@("LUPIN", 2108.0, 2402.9, 87.7273294768821, 1795.2, 117.42424242424241))
$rows += @(,@("PRESTIGE", 1597.5, 1814.0, 88.06504961411245, 1048.05, 152.42593387720052))
$rows += @(,@("YESBANK", 21.62, 24.3, 88.97119341563787, 16.02, 134.95630461922596))
$rows += @(,@("HDFCAMC", 2675.0, 2967.25, 90.15081304237931, 1781.53, 150.15183578160344))
$rows += @(,@("HEROMOTOCO", 5762.0, 6388.5, 90.19331611489395, 3344.0, 172.3086124401914))
$rows += @(,@("MPHASIS", 2788.9, 3078.4, 90.59576403326403, 2044.55, 136.4065442273361))
$rows += @(,@("SRF", 3067.7, 3325.0, 92.26165413533835, 2196.85, 139.64084939800168))
$rows += @(,@("PERSISTENT", 6257.0, 6768.8, 92.43883701690108, 4148.95, 150.80924089227395))
$rows += @(,@("PAYTM", 1297.0, 1381.8, 93.86307714575192, 651.5, 199.07904834996162))
$rows += @(,@("GMRAIRPORT", 104.2, 110.36, 94.41826748822038, 67.75, 153.8007380073801))
$rows += @(,@("MFSL", 1676.2, 1751.4, 95.70629210916981, 950.0, 176.4421052631579))
$rows += @(,@("CUMMINSIND", 4427.2, 4615.0, 95.93066088840736, 2580.0, 171.59689922480618))
$rows += @(,@("POLYCAB", 7596.0, 7903.0, 96.11539921548778, 4555.0, 166.76180021953897))
$rows += @(,@("INDUSTOWER", 421.0, 434.7, 96.8484011962273, 312.55, 134.69844824828027))
$rows += @(,@("BHEL", 288.0, 295.25, 97.54445385266723, 176.0, 163.63636363636363))
$rows += @(,@("MUTHOOTFIN", 3809.6, 3890.0, 97.9331619537275, 1965.0, 193.8727735368957))
$rows += @(,@("NMDC", 83.0, 84.64, 98.06238185255198, 59.53, 139.4254997480262))
$rows += @(,@("MARICO", 750.7, 765.3, 98.09225140467791, 577.85, 129.9126070779614))
$rows += @(,@("FEDERALBNK", 267.0, 271.1, 98.48764293618589, 172.66, 154.63917525773198))
$rows += @(,@("BHARATFORG", 1470.0, 1490.9, 98.59816218391575, 919.1, 159.93907083015992))
$rows += @(,@("APLAPOLLO", 1913.0, 1936.0, 98.81198347107438, 1272.7, 150.31036379350985))
$rows += @(,@("AUBANK", 996.4, 1007.25, 98.92280963018119, 478.35, 208.2993623915543))
$rows += @(,@("PHOENIXLTD", 1852.0, 1871.0, 98.98450026723677, 1402.5, 132.04991087344027))
$rows += @(,@("UPL", 793.05, 801.0, 99.00749063670412, 493.0, 160.86206896551724))
$rows += @(,@("ASHOKLEY", 179.0, 180.34, 99.25695907729843, 95.93, 186.59439174397997))
$rows += @(,@("IDFCFIRSTB", 85.5, 86.0, 99.4186046511628, 52.46, 162.9813191002669))
$rows += @(,@("HINDPETRO", 499.5, 500.9, 99.72050309443004, 287.55, 173.7089201877934))

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
